$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to make room for the Monk class,
# shifting the existing Rogue column (B) to column C.
$ws.Columns.Item(2).Insert()

# Populate the new Monk column (B) with the matching attribute rows:
# 1 = class name, 2 = skills, 3 = weapons, 4 = armor, 5 = equipment pack, 6 = archetype
$ws.Cells.Item(1, 2).Value() = "Monk"
$ws.Cells.Item(2, 2).Value() = "2=Acrobatics/Athletics/History/Insight/Religion/Stealth"
$ws.Cells.Item(3, 2).Value() = "Short-sword/?Simple Weapons"
$ws.Cells.Item(4, 2).Value() = "None"
$ws.Cells.Item(5, 2).Value() = "Dungeoneer's Pack/Explorer's Pack"
$ws.Cells.Item(6, 2).Value() = "3=Monastic Tradition:=classes/monk/Monastic Traditions.xlsx"
